$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: add "Description" column E ---
$ws.Range("E1").Value2 = "Description"

# --- New style for column E (center + wrap text, text format) ---
$descRange = $ws.Range("E1:E100")
$descRange.ColumnWidth = 29.54296875
$descRange.HorizontalAlignment = -4108  # xlCenter
$descRange.WrapText = $true
$descRange.NumberFormat = "@"

# --- Row data (Date, In, Out, Description) ---
$data = @(
    @{ Row = 2;  Date = 43430; In = 0.60416666666666663; Out = 0.625;               Desc = "Meeting" },
    @{ Row = 3;  Date = 43454; In = 0.61458333333333337; Out = 0.67013888888888884; Desc = "Learning dplyr/ tidytext" },
    @{ Row = 4;  Date = 43455; In = 0.44791666666666669; Out = 0.47569444444444442; Desc = "Learning dplyr" },
    @{ Row = 5;  Date = 43458; In = 0.44791666666666669; Out = 0.46875;             Desc = "Learning tidytext" },
    @{ Row = 6;  Date = 43460; In = 0.61111111111111105; Out = 0.64236111111111105; Desc = "Learning tidytext" },
    @{ Row = 7;  Date = 43462; In = 0.48958333333333331; Out = 0.51041666666666663; Desc = "Learning tidytext" },
    @{ Row = 8;  Date = 43463; In = 0.55555555555555558; Out = 0.60416666666666663; Desc = "Learning tidytext" },
    @{ Row = 9;  Date = 43469; In = 0.43402777777777773; Out = 0.45833333333333331; Desc = "Learning tidytext" },
    @{ Row = 10; Date = 43469; In = 0.71875;              Out = 0.73958333333333337; Desc = "Learning tidytext" },
    @{ Row = 11; Date = 43472; In = 0.54166666666666663; Out = 0.5625;              Desc = "Learning tidytext" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value2 = $entry.Date
    $ws.Cells.Item($r, 2).Value2 = $entry.In
    $ws.Cells.Item($r, 3).Value2 = $entry.Out
    $ws.Cells.Item($r, 5).Value2 = $entry.Desc
}

# Row 3 formula was a "non-shared" IF formula in the original file; it now
# participates in the same shared-formula group as the rest (matches the
# authored workbook's formula text).
$ws.Range("D3").Formula = '=IF(OR(ISBLANK(B3),ISBLANK(C3)), "", C3-B3)'

# --- Selection matches the authored file ---
$ws.Range("F11").Select()
